# Weekly update: a new price observation is inserted ahead of the existing
# "Femacal de La Calera - Zanahoria" series (row 337), pushing the previously
# recorded rows 337-359 down to 338-360 (exactly how Excel's row-insert
# behaves when a new data point is prepended to a time series block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 337; everything that was there (and below, down
# to the old last row 359) shifts down by one, becoming rows 338-360.
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new observation.
$ws.Range("A337").Value = 3
$ws.Range("B337").Value = "Femacal de La Calera"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = 44714
$ws.Range("E337").Value = 5
$ws.Range("F337").Value = 100114013
$ws.Range("G337").Value = "Zanahoria"
$ws.Range("H337").Value = "Sin especificar"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 410
$ws.Range("K337").Value = 7000
$ws.Range("L337").Value = 7500
$ws.Range("M337").Value = 7280
$ws.Range("N337").Value = "`$/saco 20 kilos"
$ws.Range("O337").Value = "Provincia de Quillota"
$ws.Range("P337").Value = 364
$ws.Range("Q337").Value = 20
$ws.Range("R337").Value = "Hortaliza"
